$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Paragraph "Equipamento: {{modelo}}" ( -> drop the surrounding curly quotes
# and relocate the _GoBack bookmark here, right after {{modelo}} ).
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(8)
$startPos = $p.Range.Start
$d.Range($startPos, $startPos + 1).Text = ""   # remove leading "

$p = $d.Paragraphs.Item(8)
$quotePos = $p.Range.End - 2                   # position of the trailing "
$d.Bookmarks.Add("_GoBack", $d.Range($quotePos, $quotePos))
$d.Range($quotePos, $quotePos + 1).Text = ""    # remove trailing "

# ---------------------------------------------------------------------------
# Paragraph "Número de série: {{numero_serie}}" -> drop surrounding quotes.
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(9)
$startPos = $p.Range.Start
$d.Range($startPos, $startPos + 1).Text = ""   # remove leading "

$p = $d.Paragraphs.Item(9)
$quotePos = $p.Range.End - 2                   # position of the trailing "
$d.Range($quotePos, $quotePos + 1).Text = ""    # remove trailing "

# ---------------------------------------------------------------------------
# Paragraph "Responsável pela área: {{responsavel_area}}" -> drop the
# leading quote (there never was a trailing one here).
# ---------------------------------------------------------------------------
$p = $d.Paragraphs.Item(10)
$startPos = $p.Range.Start
$d.Range($startPos, $startPos + 1).Text = ""   # remove leading "
